$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z9").Value = "2025-11-13T06:52:32.135867"
$ws.Range("Z10:Z20").Value = "2025-11-13T06:52:32.136867"
$ws.Range("Z21:Z29").Value = "2025-11-13T06:52:32.137866"
$ws.Range("Z30:Z34").Value = "2025-11-13T06:52:32.138867"
$ws.Range("Z35").Value = "2025-11-13T06:52:32.139872"
$ws.Range("Z36:Z38").Value = "2025-11-13T06:52:32.140670"
$ws.Range("Z39:Z40").Value = "2025-11-13T06:52:32.141675"
$ws.Range("Z41:Z43").Value = "2025-11-13T06:52:32.142675"
$ws.Range("Z44:Z45").Value = "2025-11-13T06:52:32.143675"
$ws.Range("Z46").Value = "2025-11-13T06:52:32.147671"
$ws.Range("Z47:Z49").Value = "2025-11-13T06:52:32.148672"
$ws.Range("Z50").Value = "2025-11-13T06:52:32.149672"
$ws.Range("Z51:Z53").Value = "2025-11-13T06:52:32.150170"
$ws.Range("Z54:Z57").Value = "2025-11-13T06:52:32.150685"
$ws.Range("Z58:Z66").Value = "2025-11-13T06:52:32.270055"
$ws.Range("Z67:Z68").Value = "2025-11-13T06:52:32.271141"
$ws.Range("Z69").Value = "2025-11-13T06:52:32.271654"
$ws.Range("Z70").Value = "2025-11-13T06:52:32.274632"
$ws.Range("Z71:Z75").Value = "2025-11-13T06:52:32.398501"
$ws.Range("Z76:Z79").Value = "2025-11-13T06:52:32.399501"
